$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
  "90×16=",
  "53×70=",
  "61×83=",
  "52×61=",
  "13×53=",
  "12×95=",
  "20×40=",
  "98×78=",
  "94×35=",
  "79×63=",
  "30×24=",
  "78×43=",
  "52×99=",
  "44×48=",
  "47×32=",
  "84×11=",
  "10×65=",
  "52×33=",
  "43×46=",
  "66×98=",
  "77×62=",
  "81×74=",
  "79×72=",
  "90×22=",
  "49×38=",
  "70×19=",
  "16×13=",
  "53×75=",
  "25×58=",
  "53×77=",
  "62×25=",
  "48×26=",
  "19×72=",
  "87×68=",
  "91×72=",
  "17×41=",
  "71×96=",
  "52×65=",
  "90×40=",
  "26×64=",
  "55×57=",
  "23×33=",
  "68×70=",
  "66×54=",
  "69×16=",
  "95×57=",
  "18×97=",
  "98×44=",
  "22×38=",
  "56×32=",
  "100×40=",
  "78×41=",
  "82×93=",
  "80×59=",
  "77×93=",
  "94×21=",
  "83×67=",
  "34×42=",
  "97×39=",
  "17×92=",
  "54×69=",
  "83×27=",
  "40×59=",
  "75×15=",
  "27×56=",
  "88×59=",
  "68×23=",
  "87×14=",
  "92×93=",
  "30×83=",
  "79×40=",
  "100×82=",
  "48×41=",
  "80×11=",
  "33×14=",
  "56×65=",
  "78×74=",
  "21×83=",
  "87×38=",
  "75×91=",
  "92×73=",
  "85×12=",
  "73×52=",
  "33×70=",
  "61×62=",
  "31×26=",
  "12×94=",
  "11×63=",
  "31×22=",
  "67×92=",
  "34×98=",
  "21×50=",
  "34×75=",
  "61×73=",
  "54×77=",
  "38×93=",
  "77×20=",
  "93×58=",
  "23×44=",
  "20×58="
)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$i]
    $i = $i + 1
  }
}
Write-Output "done: $i cells updated"